$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) cells are treated as text so values like "1.000" are not
# auto-converted to numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '28.243.96'
$ws.Range('E2').Value = '  -2.98%  '

$ws.Range('D3').Value = '1.918.44'
$ws.Range('E3').Value = '  -3.69%  '

$ws.Range('D4').Value = '1.000'
$ws.Range('E4').Value = '  -1.35%  '

$ws.Range('D5').Value = '327.99'
$ws.Range('E5').Value = '  -0.76%  '

$ws.Range('D6').Value = '0.9993'
$ws.Range('E6').Value = '  -1.18%  '

$ws.Range('D7').Value = '0.4689'
$ws.Range('E7').Value = '  -5.69%  '

$ws.Range('D8').Value = '0.4031'
$ws.Range('E8').Value = '  -3.85%  '

$ws.Range('D9').Value = '53.10'
$ws.Range('E9').Value = '  -2.93%  '

$ws.Range('D10').Value = '0.08439'
$ws.Range('E10').Value = '  -5.26%  '

$ws.Range('D11').Value = '1.047'
$ws.Range('E11').Value = '  -4.40%  '

$ws.Range('D12').Value = '22.24'
$ws.Range('E12').Value = '  -4.29%  '

$ws.Range('D13').Value = '1.947.38'
$ws.Range('E13').Value = '  -2.10%  '

$ws.Range('D14').Value = '7.451'
$ws.Range('E14').Value = '  -6.93%  '

$ws.Range('D15').Value = '6.084'
$ws.Range('E15').Value = '  -5.48%  '

$ws.Range('D16').Value = '1.002'
$ws.Range('E16').Value = '  -1.28%  '

$ws.Range('D17').Value = '89.83'
$ws.Range('E17').Value = '  -2.83%  '

$ws.Range('D18').Value = '0.00001066'
$ws.Range('E18').Value = '  -3.61%  '

$ws.Range('D19').Value = '0.06594'
$ws.Range('E19').Value = '  -2.51%  '

$ws.Range('D20').Value = '18.02'
$ws.Range('E20').Value = '  -7.99%  '

$ws.Range('D21').Value = '0.9996'
$ws.Range('E21').Value = '  -1.23%  '

$ws.Range('D22').Value = '5.730'
$ws.Range('E22').Value = '  -4.08%  '

$ws.Range('D23').Value = '28.275.00'
$ws.Range('E23').Value = '  -2.89%  '

$ws.Range('D24').Value = '11.30'
$ws.Range('E24').Value = '  -5.77%  '

$ws.Range('D25').Value = '2.277'
$ws.Range('E25').Value = '  -0.55%  '

$ws.Range('D26').Value = '2.166.47'
$ws.Range('E26').Value = '  -2.84%  '

$ws.Range('D27').Value = '153.23'
$ws.Range('E27').Value = '  -2.51%  '

$ws.Range('D28').Value = '20.03'
$ws.Range('E28').Value = '  -3.83%  '

$ws.Range('D29').Value = '5.748'
$ws.Range('E29').Value = '  -8.63%  '

$ws.Range('D30').Value = '2.125'
$ws.Range('E30').Value = '  -5.99%  '

$ws.Range('D31').Value = '123.66'
$ws.Range('E31').Value = '  -2.94%  '

$ws.Range('D32').Value = '0.9773'
$ws.Range('E32').Value = '  -6.66%  '

$ws.Range('D33').Value = '0.09637'
$ws.Range('E33').Value = '  -2.10%  '

$ws.Range('D34').Value = '1.457'
$ws.Range('E34').Value = '  -5.05%  '

$ws.Range('B35').Value = 'Filecoin'
$ws.Range('C35').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D35').Value = '5.552'
$ws.Range('E35').Value = '  -4.70%  '

$ws.Range('B36').Value = 'HuobiToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D36').Value = '3.630'
$ws.Range('E36').Value = '  -3.07%  '

$ws.Range('D37').Value = '8.814'
$ws.Range('E37').Value = '  -2.85%  '

$ws.Range('D38').Value = '0.02301'
$ws.Range('E38').Value = '  -5.04%  '

$ws.Range('E39').Value = '  -4.56%  '

$ws.Range('D40').Value = '0.06171'
$ws.Range('E40').Value = '  -3.54%  '

$ws.Range('B41').Value = 'TheSandbox'
$ws.Range('C41').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D41').Value = '0.6148'
$ws.Range('E41').Value = '  -5.15%  '

$ws.Range('B42').Value = 'Aptos'
$ws.Range('C42').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D42').Value = '11.07'
$ws.Range('E42').Value = '  -3.58%  '

$ws.Range('D43').Value = '0.9992'
$ws.Range('E43').Value = '  -1.06%  '

$ws.Range('D44').Value = '0.1906'

$ws.Range('D45').Value = '1.301'
$ws.Range('E45').Value = '  -3.35%  '

$ws.Range('B46').Value = 'EnergySwap'
$ws.Range('C46').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D46').Value = '12.84'
$ws.Range('E46').Value = '  -4.60%  '

$ws.Range('B47').Value = 'Decentraland'
$ws.Range('C47').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D47').Value = '0.5844'
$ws.Range('E47').Value = '  -5.93%  '

$ws.Range('D48').Value = '2.028'
$ws.Range('E48').Value = '  -7.19%  '

$ws.Range('D49').Value = '3.427'
$ws.Range('E49').Value = '  -1.81%  '

$ws.Range('D50').Value = '0.06900'
$ws.Range('E50').Value = '  -1.23%  '

$ws.Range('D51').Value = '110.23'
$ws.Range('E51').Value = '  -2.55%  '
